$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly-completed match row (row 45): actual/predicted winners,
# man of the match, and amount won. Columns B:E are shared-string text
# cells that previously didn't exist, so give them the same centered style
# used across the rest of the table before writing values.
$ws.Range("B45:G45").HorizontalAlignment = -4108

$ws.Range("B45").Value = "MI"
$ws.Range("C45").Value = "Rohit"
$ws.Range("D45").Value = "MI"
$ws.Range("E45").Value = "QDK"
$ws.Range("F45").Value = 3
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 90

# The hidden AutoFilter defined name needs to grow to cover the newly
# completed row.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$J`$45"
    }
}

# Update the saved selection/scroll position to match the author's last
# cursor position when the file was saved.
$ws.Range("E6").Select()
